$wb = $excel.ActiveWorkbook

# --- Add the new "Variable Definitions" sheet after the last existing sheet (HurdleRates) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Variable Definitions"

# --- Populate the variable-definitions table (order matches shared-string insertion order) ---
$ws.Range("A1").Value = "Variable Name"
$ws.Range("B1").Value = "Description"

$ws.Range("A2").Value = "`$/theta`$"
$ws.Range("A3").Value = "`$n`$"
$ws.Range("A4").Value = "`$p`$"

$ws.Range("B4").Value = "Per-attempt probability of success"
$ws.Range("B3").Value = "Number of attempts"
$ws.Range("B2").Value = "Target probability of success"

$ws.Range("A5").Value = "`$X`$"
$ws.Range("B5").Value = "Present value of pull size"

# --- Set page orientation to portrait on the HurdleRates sheet ---
$wsHurdle = $wb.Worksheets.Item("HurdleRates")
$wsHurdle.PageSetup.Orientation = 1

# --- Update selections / active sheet to match the final view state ---
$wsHurdle.Activate()
$wsHurdle.Range("A15").Select()

$ws.Activate()
$ws.Range("E30").Select()
